$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "List" column (E) values for rows 2-4 to the new semicolon list
$ws.Range("E2").Value = "ayush;deepanshu;aryan;sumit"
$ws.Range("E3").Value = "ayush;deepanshu;aryan;sumit"
$ws.Range("E4").Value = "ayush;deepanshu;aryan;sumit"

# Add new "Quantity" column (F)
$ws.Range("F1").Value = "Quantity"
$ws.Range("F2").Value = 10
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 7

# Update the active cell selection to match the target state
$ws.Range("K5").Select()
